# Apply row-content permutation changes to rows 7-21 of the 'Artfynd' sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 1).Value = 111936855
$ws.Cells.Item(7, 2).Value = 78579
$ws.Cells.Item(7, 4).Value = 'NT'
$ws.Cells.Item(7, 5).Value = 2081
$ws.Cells.Item(7, 6).Value = 'Skrovellav'
$ws.Cells.Item(7, 7).Value = 'Lobaria scrobiculata'
$ws.Cells.Item(7, 8).Value = '(Scop.) DC.'
$ws.Cells.Item(7, 17).Value = 449178.0024977843
$ws.Cells.Item(7, 18).Value = 7087509.952929454
$ws.Cells.Item(8, 1).Value = 111936790
$ws.Cells.Item(8, 2).Value = 90087
$ws.Cells.Item(8, 4).Value = 'LC'
$ws.Cells.Item(8, 5).Value = 3298
$ws.Cells.Item(8, 6).Value = 'Trådticka'
$ws.Cells.Item(8, 7).Value = 'Climacocystis borealis'
$ws.Cells.Item(8, 8).Value = '(Fr.) Kotl. & Pouzar'
$ws.Cells.Item(8, 17).Value = 448943.8940418276
$ws.Cells.Item(8, 18).Value = 7087698.235520792
$ws.Cells.Item(9, 1).Value = 111936856
$ws.Cells.Item(9, 2).Value = 78579
$ws.Cells.Item(9, 5).Value = 2081
$ws.Cells.Item(9, 6).Value = 'Skrovellav'
$ws.Cells.Item(9, 7).Value = 'Lobaria scrobiculata'
$ws.Cells.Item(9, 8).Value = '(Scop.) DC.'
$ws.Cells.Item(9, 17).Value = 449235.6140813087
$ws.Cells.Item(9, 18).Value = 7087446.93781954
$ws.Cells.Item(10, 1).Value = 111936786
$ws.Cells.Item(10, 2).Value = 89405
$ws.Cells.Item(10, 5).Value = 1202
$ws.Cells.Item(10, 6).Value = 'Ullticka'
$ws.Cells.Item(10, 7).Value = 'Phellinidium ferrugineofuscum'
$ws.Cells.Item(10, 8).Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Cells.Item(10, 17).Value = 449055.4908092887
$ws.Cells.Item(10, 18).Value = 7087729.435905921
$ws.Cells.Item(10, 29).Value = $null
$ws.Cells.Item(11, 1).Value = 111936802
$ws.Cells.Item(11, 2).Value = 56398
$ws.Cells.Item(11, 5).Value = 100109
$ws.Cells.Item(11, 6).Value = 'Tretåig hackspett'
$ws.Cells.Item(11, 7).Value = 'Picoides tridactylus'
$ws.Cells.Item(11, 8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(11, 17).Value = 448825.396353531
$ws.Cells.Item(11, 18).Value = 7087649.541088912
$ws.Cells.Item(11, 29).Value = 'ringhack äldre'
$ws.Cells.Item(12, 1).Value = 111936874
$ws.Cells.Item(12, 2).Value = 89423
$ws.Cells.Item(12, 5).Value = 5432
$ws.Cells.Item(12, 6).Value = 'Granticka'
$ws.Cells.Item(12, 7).Value = 'Porodaedalea chrysoloma'
$ws.Cells.Item(12, 8).Value = '(Fr.) Fiasson & Niemelä'
$ws.Cells.Item(12, 17).Value = 449309.3031681653
$ws.Cells.Item(12, 18).Value = 7087423.7290034
$ws.Cells.Item(13, 1).Value = 111936877
$ws.Cells.Item(13, 2).Value = 89423
$ws.Cells.Item(13, 5).Value = 5432
$ws.Cells.Item(13, 6).Value = 'Granticka'
$ws.Cells.Item(13, 7).Value = 'Porodaedalea chrysoloma'
$ws.Cells.Item(13, 8).Value = '(Fr.) Fiasson & Niemelä'
$ws.Cells.Item(13, 17).Value = 449302.006836799
$ws.Cells.Item(13, 18).Value = 7087517.627985355
$ws.Cells.Item(13, 29).Value = $null
$ws.Cells.Item(14, 1).Value = 111936876
$ws.Cells.Item(14, 2).Value = 89423
$ws.Cells.Item(14, 5).Value = 5432
$ws.Cells.Item(14, 6).Value = 'Granticka'
$ws.Cells.Item(14, 7).Value = 'Porodaedalea chrysoloma'
$ws.Cells.Item(14, 8).Value = '(Fr.) Fiasson & Niemelä'
$ws.Cells.Item(14, 17).Value = 449317.0489210376
$ws.Cells.Item(14, 18).Value = 7087520.906595955
$ws.Cells.Item(15, 1).Value = 111936800
$ws.Cells.Item(15, 17).Value = 449177.8905366624
$ws.Cells.Item(15, 18).Value = 7087530.207590466
$ws.Cells.Item(16, 1).Value = 111936879
$ws.Cells.Item(16, 17).Value = 449280.5263497296
$ws.Cells.Item(16, 18).Value = 7087552.317575688
$ws.Cells.Item(17, 1).Value = 111936791
$ws.Cells.Item(17, 2).Value = 90087
$ws.Cells.Item(17, 4).Value = 'LC'
$ws.Cells.Item(17, 5).Value = 3298
$ws.Cells.Item(17, 6).Value = 'Trådticka'
$ws.Cells.Item(17, 7).Value = 'Climacocystis borealis'
$ws.Cells.Item(17, 8).Value = '(Fr.) Kotl. & Pouzar'
$ws.Cells.Item(17, 17).Value = 449279.9868849564
$ws.Cells.Item(17, 18).Value = 7087437.412661138
$ws.Cells.Item(18, 1).Value = 111936801
$ws.Cells.Item(18, 2).Value = 56398
$ws.Cells.Item(18, 5).Value = 100109
$ws.Cells.Item(18, 6).Value = 'Tretåig hackspett'
$ws.Cells.Item(18, 7).Value = 'Picoides tridactylus'
$ws.Cells.Item(18, 8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(18, 17).Value = 448837.1356602872
$ws.Cells.Item(18, 18).Value = 7087666.519360777
$ws.Cells.Item(18, 29).Value = 'ringhack äldre'
$ws.Cells.Item(19, 1).Value = 111936873
$ws.Cells.Item(19, 2).Value = 89423
$ws.Cells.Item(19, 5).Value = 5432
$ws.Cells.Item(19, 6).Value = 'Granticka'
$ws.Cells.Item(19, 7).Value = 'Porodaedalea chrysoloma'
$ws.Cells.Item(19, 8).Value = '(Fr.) Fiasson & Niemelä'
$ws.Cells.Item(19, 17).Value = 449170.9149442808
$ws.Cells.Item(19, 18).Value = 7087507.866619493
$ws.Cells.Item(20, 1).Value = 111936872
$ws.Cells.Item(20, 2).Value = 89423
$ws.Cells.Item(20, 4).Value = 'NT'
$ws.Cells.Item(20, 5).Value = 5432
$ws.Cells.Item(20, 6).Value = 'Granticka'
$ws.Cells.Item(20, 7).Value = 'Porodaedalea chrysoloma'
$ws.Cells.Item(20, 8).Value = '(Fr.) Fiasson & Niemelä'
$ws.Cells.Item(20, 17).Value = 449151.0056780232
$ws.Cells.Item(20, 18).Value = 7087530.644260203
$ws.Cells.Item(21, 1).Value = 111936895
$ws.Cells.Item(21, 2).Value = 85715
$ws.Cells.Item(21, 5).Value = 510
$ws.Cells.Item(21, 6).Value = 'Doftskinn'
$ws.Cells.Item(21, 7).Value = 'Cystostereum murrayi'
$ws.Cells.Item(21, 8).Value = '(Berk. & M.A. Curtis.) Pouzar'
$ws.Cells.Item(21, 17).Value = 448924.849606293
$ws.Cells.Item(21, 18).Value = 7087773.829047815
